$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("AC3").Value = "x"
$r = $ws.Range("AC3")
$r.Interior.TintAndShade = -0.14999847407452621
$r.Interior.ThemeColor = 2
